$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'', RandomUnderSampler(random_state=42)),
                (''model'',
                 BaggingClassifier(estimator=DecisionTreeClassifier(criterion=''entropy'',
                                                                    max_depth=1,
                                                                    max_features=''sqrt'',
                                                                    min_samples_leaf=4,
                                                                    min_samples_split=4,
                                                                    random_state=42),
                                   n_estimators=5, random_state=42))])'
$ws.Range("B2").Value = 0.6696266968325791
$ws.Range("C2").Value = '{''selector'': RandomUnderSampler(random_state=42), ''scaler'': MinMaxScaler(), ''model__n_estimators'': 5, ''model__estimator__min_samples_split'': 4, ''model__estimator__min_samples_leaf'': 4, ''model__estimator__max_features'': ''sqrt'', ''model__estimator__max_depth'': 1, ''model__estimator__criterion'': ''entropy'', ''model__estimator__class_weight'': None}'
$ws.Range("F2").Value = '[1 1 1 1 1 1 1 1 1 0 0 1]'
$ws.Range("H2").Value = 0.8319604881314266
$ws.Range("I2").Value = 0.02812033012710329
$ws.Range("J2").Value = 0.5065359610977257
$ws.Range("K2").Value = 0.1384243984303357

# Row 3
$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'', RandomUnderSampler(random_state=42)),
                (''model'',
                 BaggingClassifier(estimator=DecisionTreeClassifier(class_weight=''balanced'',
                                                                    criterion=''entropy'',
                                                                    max_depth=1,
                                                                    min_samples_leaf=6,
                                                                    min_samples_split=6,
                                                                    random_state=42),
                                   random_state=42))])'
$ws.Range("B3").Value = 0.7249768368886016
$ws.Range("C3").Value = '{''selector'': RandomUnderSampler(random_state=42), ''scaler'': MinMaxScaler(), ''model__n_estimators'': 10, ''model__estimator__min_samples_split'': 6, ''model__estimator__min_samples_leaf'': 6, ''model__estimator__max_features'': None, ''model__estimator__max_depth'': 1, ''model__estimator__criterion'': ''entropy'', ''model__estimator__class_weight'': ''balanced''}'
$ws.Range("D3").Value = 0.3076923076923077
$ws.Range("F3").Value = '[0 0 0 1 1 1 0 1 0 1 0 1]'
$ws.Range("H3").Value = 0.8251475159617845
$ws.Range("I3").Value = 0.03403574868529025
$ws.Range("J3").Value = 0.6429551498175028
$ws.Range("K3").Value = 0.1179955866028993

# Row 4
$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'', RandomUnderSampler(random_state=42)),
                (''model'',
                 BaggingClassifier(estimator=DecisionTreeClassifier(class_weight=''balanced'',
                                                                    criterion=''entropy'',
                                                                    max_depth=2,
                                                                    max_features=''sqrt'',
                                                                    min_samples_leaf=6,
                                                                    min_samples_split=3,
                                                                    random_state=42),
                                   random_state=42))])'
$ws.Range("B4").Value = 0.7380494505494506
$ws.Range("C4").Value = '{''selector'': RandomUnderSampler(random_state=42), ''scaler'': MinMaxScaler(), ''model__n_estimators'': 10, ''model__estimator__min_samples_split'': 3, ''model__estimator__min_samples_leaf'': 6, ''model__estimator__max_features'': ''sqrt'', ''model__estimator__max_depth'': 2, ''model__estimator__criterion'': ''entropy'', ''model__estimator__class_weight'': ''balanced''}'
$ws.Range("D4").Value = 0.823529411764706
$ws.Range("F4").Value = '[0 1 1 1 1 1 0 1 1 1 1 0]'
$ws.Range("H4").Value = 0.8298205206197684
$ws.Range("I4").Value = 0.02939855730219668
$ws.Range("J4").Value = 0.6078983170424347
$ws.Range("K4").Value = 0.1278670651110609

# Row 5
$ws.Range("A5").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'', RandomUnderSampler(random_state=42)),
                (''model'',
                 BaggingClassifier(estimator=DecisionTreeClassifier(class_weight=''balanced'',
                                                                    criterion=''entropy'',
                                                                    max_depth=2,
                                                                    max_features=''log2'',
                                                                    min_samples_split=5,
                                                                    random_state=42),
                                   n_estimators=5, random_state=42))])'
$ws.Range("B5").Value = 0.7174768368886015
$ws.Range("C5").Value = '{''selector'': RandomUnderSampler(random_state=42), ''scaler'': MinMaxScaler(), ''model__n_estimators'': 5, ''model__estimator__min_samples_split'': 5, ''model__estimator__min_samples_leaf'': 1, ''model__estimator__max_features'': ''log2'', ''model__estimator__max_depth'': 2, ''model__estimator__criterion'': ''entropy'', ''model__estimator__class_weight'': ''balanced''}'
$ws.Range("D5").Value = 0.5333333333333333
$ws.Range("F5").Value = '[1 1 1 1 0 0 1 1 1 1 1 1]'
$ws.Range("H5").Value = 0.8104996944284236
$ws.Range("I5").Value = 0.02704202896776954
$ws.Range("J5").Value = 0.623978951767187
$ws.Range("K5").Value = 0.1452569815811669
